# Refresh the cryptos price table (columns B-E, rows 2-51) to match the
# latest scrape pulled in by the GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells hold digit-only text (e.g. "1.00", "0.210"); without
# forcing Text format first, Excel would coerce them to numbers and drop
# the significant trailing zeros that the source formatting relies on.
$textRanges = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D12", "D13", "D14", "D18", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($r in $textRanges) { $ws.Range($r).NumberFormat = "@" }

$ws.Range('D2').Value = '96.327.90'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.643.41'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('D4').Value = '2.61'
$ws.Range('E4').Value = '  +36.17%  '
$ws.Range('D5').Value = '1.00'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = '226.21'
$ws.Range('E6').Value = '  -4.85%  '
$ws.Range('D7').Value = '640.87'
$ws.Range('E7').Value = '  -2.97%  '
$ws.Range('D8').Value = '0.422'
$ws.Range('E8').Value = '  -0.86%  '
$ws.Range('D9').Value = '1.13'
$ws.Range('E9').Value = '  +5.70%  '
$ws.Range('D10').Value = '1.00'
$ws.Range('E10').Value = '  +0.09%  '
$ws.Range('D11').Value = '3.644.71'
$ws.Range('E11').Value = '  -1.55%  '
$ws.Range('D12').Value = '48.52'
$ws.Range('E12').Value = '  +9.07%  '
$ws.Range('D13').Value = '0.210'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '0.0000292'
$ws.Range('E14').Value = '  -9.57%  '
$ws.Range('E15').Value = '  -3.56%  '
$ws.Range('D16').Value = '4.327.95'
$ws.Range('E16').Value = '  -1.62%  '
$ws.Range('D17').Value = '95.823.93'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').Value = '21.46'
$ws.Range('E18').Value = '  +14.27%  '
$ws.Range('D19').Value = '8.86'
$ws.Range('E19').Value = '  -2.70%  '
$ws.Range('D20').Value = '14.18'
$ws.Range('E20').Value = '  +8.41%  '
$ws.Range('D21').Value = '3.649.26'
$ws.Range('E21').Value = '  -1.56%  '
$ws.Range('D22').Value = '0.537'
$ws.Range('E22').Value = '  +6.30%  '
$ws.Range('D23').Value = '0.260'
$ws.Range('E23').Value = '  +33.38%  '
$ws.Range('D24').Value = '515.43'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').Value = '3.27'
$ws.Range('E25').Value = '  -5.21%  '
$ws.Range('D26').Value = '122.10'
$ws.Range('E26').Value = '  +18.24%  '
$ws.Range('D27').Value = '0.0000202'
$ws.Range('E27').Value = '  -7.49%  '
$ws.Range('D28').Value = '6.81'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('D29').Value = '12.89'
$ws.Range('E29').Value = '  -4.62%  '
$ws.Range('D30').Value = '13.17'
$ws.Range('E30').Value = '  +1.51%  '
$ws.Range('D31').Value = '2.99'
$ws.Range('E31').Value = '  -2.10%  '
$ws.Range('E32').Value = '  +0.02%  '
$ws.Range('E33').Value = '  -4.89%  '
$ws.Range('D34').Value = '33.03'
$ws.Range('E34').Value = '  +2.35%  '
$ws.Range('D35').Value = '0.620'
$ws.Range('E35').Value = '  +4.02%  '
$ws.Range('D36').Value = '0.996'
$ws.Range('E36').Value = '  -0.69%  '
$ws.Range('D37').Value = '1.77'
$ws.Range('E37').Value = '  -4.75%  '
$ws.Range('B38').Value = 'USDe'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.02%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '600.35'
$ws.Range('E39').Value = '  -8.12%  '
$ws.Range('D40').Value = '8.46'
$ws.Range('E40').Value = '  -4.41%  '
$ws.Range('D41').Value = '43.46'
$ws.Range('E41').Value = '  +7.85%  '
$ws.Range('D42').Value = '7.10'
$ws.Range('E42').Value = '  +3.08%  '
$ws.Range('D43').Value = '0.490'
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '0.160'
$ws.Range('E44').Value = '  -5.07%  '
$ws.Range('D45').Value = '0.0497'
$ws.Range('E45').Value = '  +7.32%  '
$ws.Range('D46').Value = '0.960'
$ws.Range('E46').Value = '  -1.33%  '
$ws.Range('D47').Value = '1.95'
$ws.Range('E47').Value = '  -5.28%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '228.59'
$ws.Range('E48').Value = '  +11.15%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.29'
$ws.Range('E49').Value = '  -5.99%  '
$ws.Range('D50').Value = '8.77'
$ws.Range('E50').Value = '  -0.39%  '
$ws.Range('D51').Value = '23.57'
